# Applies updated market-price / profit figures to the Leve profit
# tables on each job sheet, as produced by the scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1689.0769
$ws.Range("I2").Value = 453.75
$ws.Range("K2").Value = 453.75
$ws.Range("M2").Value = -340.75
$ws.Range("H5").Value = 102
$ws.Range("I5").Value = 49.25
$ws.Range("K5").Value = 49.25
$ws.Range("M5").Value = 65.75
$ws.Range("H12").Value = 425.91666
$ws.Range("I12").Value = 442.2857
$ws.Range("J12").Value = 403
$ws.Range("K12").Value = 442.2857
$ws.Range("L12").Value = 403
$ws.Range("M12").Value = -272.2857
$ws.Range("N12").Value = -743
$ws.Range("H33").Value = 4283.2
$ws.Range("I33").Value = 5394.2104
$ws.Range("J33").Value = 765
$ws.Range("K33").Value = 5394.2104
$ws.Range("L33").Value = 765
$ws.Range("M33").Value = -5165.2104
$ws.Range("N33").Value = -1223
$ws.Range("H40").Value = 2763.2727
$ws.Range("I40").Value = 1517.6471
$ws.Range("J40").Value = 6998.4
$ws.Range("K40").Value = 1517.6471
$ws.Range("L40").Value = 6998.4
$ws.Range("M40").Value = -1342.6471
$ws.Range("N40").Value = -7348.4
$ws.Range("H55").Value = 672
$ws.Range("I55").Value = 73.666664
$ws.Range("K55").Value = 73.666664
$ws.Range("M55").Value = 140.333336
$ws.Range("H94").Value = 4754.2856
$ws.Range("I94").Value = 4380
$ws.Range("K94").Value = 4380
$ws.Range("M94").Value = -3929
$ws.Range("H98").Value = 2504.606
$ws.Range("I98").Value = 2457.9062
$ws.Range("K98").Value = 2457.9062
$ws.Range("M98").Value = -959.9061999999999
$ws.Range("H99").Value = 5340
$ws.Range("I99").Value = 181
$ws.Range("K99").Value = 543
$ws.Range("M99").Value = 955
$ws.Range("H122").Value = 2504.606
$ws.Range("I122").Value = 2457.9062
$ws.Range("K122").Value = 7373.7186
$ws.Range("M122").Value = -4923.7186
$ws.Range("H138").Value = 15642.944
$ws.Range("I138").Value = 18154.316
$ws.Range("J138").Value = 9678.4375
$ws.Range("K138").Value = 54462.948
$ws.Range("L138").Value = 29035.3125
$ws.Range("M138").Value = -49322.948
$ws.Range("N138").Value = -39315.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 271.53845
$ws.Range("J5").Value = 198.16667
$ws.Range("L5").Value = 198.16667
$ws.Range("N5").Value = -422.16667
$ws.Range("H26").Value = 1007
$ws.Range("I26").Value = 1007
$ws.Range("K26").Value = 1007
$ws.Range("M26").Value = -677
$ws.Range("H32").Value = 23718.809
$ws.Range("I32").Value = 24910.521
$ws.Range("K32").Value = 24910.521
$ws.Range("M32").Value = -24623.521
$ws.Range("H45").Value = 3679.25
$ws.Range("I45").Value = 4499.5
$ws.Range("K45").Value = 4499.5
$ws.Range("M45").Value = -4122.5
$ws.Range("H61").Value = 9752.75
$ws.Range("I61").Value = 9505.5
$ws.Range("K61").Value = 9505.5
$ws.Range("M61").Value = -9293.5
$ws.Range("H97").Value = 14773.667
$ws.Range("I97").Value = 34600
$ws.Range("J97").Value = 4860.5
$ws.Range("K97").Value = 34600
$ws.Range("L97").Value = 4860.5
$ws.Range("M97").Value = -34104
$ws.Range("N97").Value = -5852.5
$ws.Range("H110").Value = 3237
$ws.Range("I110").Value = 2117.4546
$ws.Range("K110").Value = 2117.4546
$ws.Range("M110").Value = -72.45460000000003
$ws.Range("H122").Value = 1130.1818
$ws.Range("I122").Value = 1108
$ws.Range("J122").Value = 1474
$ws.Range("K122").Value = 3324
$ws.Range("L122").Value = 4422
$ws.Range("M122").Value = -874
$ws.Range("N122").Value = -9322
$ws.Range("H136").Value = 9752.75
$ws.Range("I136").Value = 9505.5
$ws.Range("K136").Value = 28516.5
$ws.Range("M136").Value = -25966.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 271.53845
$ws.Range("J4").Value = 198.16667
$ws.Range("L4").Value = 198.16667
$ws.Range("N4").Value = -428.16667
$ws.Range("H5").Value = 2135.2856
$ws.Range("I5").Value = 3511.75
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 3511.75
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -3398.75
$ws.Range("N5").Value = -526
$ws.Range("H7").Value = 8334002
$ws.Range("I7").Value = 1003
$ws.Range("K7").Value = 1003
$ws.Range("M7").Value = -890
$ws.Range("H80").Value = 346.2857
$ws.Range("I80").Value = 346.66666
$ws.Range("K80").Value = 346.66666
$ws.Range("M80").Value = 651.33334
$ws.Range("H83").Value = 346.2857
$ws.Range("I83").Value = 346.66666
$ws.Range("K83").Value = 1733.3333
$ws.Range("M83").Value = 3258.6667
$ws.Range("H94").Value = 4656.4546
$ws.Range("I94").Value = 3979.4
$ws.Range("J94").Value = 5220.6665
$ws.Range("K94").Value = 3979.4
$ws.Range("L94").Value = 5220.6665
$ws.Range("M94").Value = -3528.4
$ws.Range("N94").Value = -6122.6665
$ws.Range("H134").Value = 2932.739
$ws.Range("I134").Value = 2932.739
$ws.Range("K134").Value = 8798.217000000001
$ws.Range("M134").Value = -6263.217000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 80762.08
$ws.Range("I58").Value = 103626.5
$ws.Range("J58").Value = 4547.3335
$ws.Range("K58").Value = 103626.5
$ws.Range("L58").Value = 4547.3335
$ws.Range("M58").Value = -103423.5
$ws.Range("N58").Value = -4953.3335
$ws.Range("H132").Value = 2497.15
$ws.Range("I132").Value = 2209.6956
$ws.Range("K132").Value = 6629.0868
$ws.Range("M132").Value = -4099.0868
$ws.Range("H134").Value = 102908.1
$ws.Range("I134").Value = 251499.75
$ws.Range("K134").Value = 754499.25
$ws.Range("M134").Value = -751964.25
$ws.Range("H136").Value = 80762.08
$ws.Range("I136").Value = 103626.5
$ws.Range("J136").Value = 4547.3335
$ws.Range("K136").Value = 310879.5
$ws.Range("L136").Value = 13642.0005
$ws.Range("M136").Value = -308329.5
$ws.Range("N136").Value = -18742.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 548
$ws.Range("I38").Value = 465
$ws.Range("K38").Value = 1395
$ws.Range("M38").Value = -1048
$ws.Range("H107").Value = 420.4
$ws.Range("J107").Value = 420.4
$ws.Range("L107").Value = 1261.2
$ws.Range("N107").Value = -5101.2
$ws.Range("H122").Value = 623.5454999999999
$ws.Range("J122").Value = 743.3333
$ws.Range("L122").Value = 6689.9997
$ws.Range("N122").Value = -11589.9997
$ws.Range("H140").Value = 2653.7273
$ws.Range("I140").Value = 2119.1
$ws.Range("K140").Value = 6357.299999999999
$ws.Range("M140").Value = -1177.299999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 167067.67
$ws.Range("I3").Value = 301.5
$ws.Range("J3").Value = 250450.75
$ws.Range("K3").Value = 301.5
$ws.Range("L3").Value = 250450.75
$ws.Range("M3").Value = -185.5
$ws.Range("N3").Value = -250682.75
$ws.Range("H4").Value = 99
$ws.Range("I4").Value = 99
$ws.Range("K4").Value = 99
$ws.Range("M4").Value = 13
$ws.Range("H122").Value = 4157.625
$ws.Range("I122").Value = 2487.125
$ws.Range("J122").Value = 7498.625
$ws.Range("K122").Value = 7461.375
$ws.Range("L122").Value = 22495.875
$ws.Range("M122").Value = -5011.375
$ws.Range("N122").Value = -27395.875
$ws.Range("H126").Value = 5410.6665
$ws.Range("I126").Value = 4717.273
$ws.Range("J126").Value = 7317.5
$ws.Range("K126").Value = 14151.819
$ws.Range("L126").Value = 21952.5
$ws.Range("M126").Value = -11681.819
$ws.Range("N126").Value = -26892.5
$ws.Range("H132").Value = 147311.42
$ws.Range("I132").Value = 204156.2
$ws.Range("J132").Value = 5199.5
$ws.Range("K132").Value = 612468.6000000001
$ws.Range("L132").Value = 15598.5
$ws.Range("M132").Value = -609938.6000000001
$ws.Range("N132").Value = -20658.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4835.5
$ws.Range("I61").Value = 1367.5264
$ws.Range("K61").Value = 1367.5264
$ws.Range("M61").Value = -1165.5264
$ws.Range("H100").Value = 3724.375
$ws.Range("I100").Value = 1999.6666
$ws.Range("J100").Value = 4759.2
$ws.Range("K100").Value = 1999.6666
$ws.Range("L100").Value = 4759.2
$ws.Range("M100").Value = -1458.6666
$ws.Range("N100").Value = -5841.2
$ws.Range("H113").Value = 4835.5
$ws.Range("I113").Value = 1367.5264
$ws.Range("K113").Value = 1367.5264
$ws.Range("M113").Value = 802.4736
$ws.Range("H122").Value = 4496.7393
$ws.Range("I122").Value = 3498.7778
$ws.Range("J122").Value = 5138.2856
$ws.Range("K122").Value = 10496.3334
$ws.Range("L122").Value = 15414.8568
$ws.Range("M122").Value = -8046.3334
$ws.Range("N122").Value = -20314.8568
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H113").Value = 1633.1428
$ws.Range("I113").Value = 383
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 1149
$ws.Range("L113").Value = 9900
$ws.Range("M113").Value = 1021
$ws.Range("N113").Value = -14240
$ws.Range("H122").Value = 9795.666999999999
$ws.Range("I122").Value = 12504.223
$ws.Range("J122").Value = 1670
$ws.Range("K122").Value = 37512.669
$ws.Range("L122").Value = 5010
$ws.Range("M122").Value = -35062.669
$ws.Range("N122").Value = -9910
$ws.Range("H126").Value = 43628.48
$ws.Range("I126").Value = 59017.445
$ws.Range("K126").Value = 59017.445
$ws.Range("M126").Value = -174582.335
